# Generate Report for Handback
#
# For each language sheet (zh-cn, de-de):
#   - Status (column B, row 2) moves from "Ready for handoff" to
#     "Handed back: in sync with en-US".
#   - The "Latest Target File" (E2) and "Latest Handback File" (F2) columns
#     get populated with hyperlinked filenames, mirroring the existing
#     "Source File Name" (A2) / "Latest Handoff File" (C2) hyperlinks.
#   - "Latest Handback DateTime" (G2) is stamped with the handback time.

$wb = $excel.ActiveWorkbook

# Color used by the workbook's existing custom "HyperLink" cell style
# (rgb="FF6495ED"). Range.Font.Color is read back out as BGR, so the
# integer we need to assign is the byte-reversed value of 0x6495ED.
$hyperlinkColor = 15570276

function Set-HandbackRow($SheetName, $StatusText, $TargetFileName, $TargetFileUrl, $HandbackFileName, $HandbackFileUrl, $HandbackDateTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    # Status: "Ready for handoff" -> "Handed back: in sync with en-US"
    $ws.Range("B2").Value = $StatusText

    # Latest Target File (E2) - same target as the Source File Name (A2)
    $ws.Range("E2").Value = $TargetFileName
    $ws.Hyperlinks.Add($ws.Range("E2"), $TargetFileUrl, "", "", $TargetFileName)
    $ws.Range("E2").Font.Underline = $true
    $ws.Range("E2").Font.Color = $hyperlinkColor

    # Latest Handback File (F2) - same target as the Latest Handoff File (C2)
    $ws.Range("F2").Value = $HandbackFileName
    $ws.Hyperlinks.Add($ws.Range("F2"), $HandbackFileUrl, "", "", $HandbackFileName)
    $ws.Range("F2").Font.Underline = $true
    $ws.Range("F2").Font.Color = $hyperlinkColor

    # Latest Handback DateTime (G2)
    $ws.Range("G2").Value = $HandbackDateTime
}

Set-HandbackRow "zh-cn" `
    "Handed back: in sync with en-US" `
    "c69f6660-4d93-4805-9fc9-43f02af01ecc.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/55e2dd23eb0bd25de4cc8458b6ae18eef6df06dd/e2e/c69f6660-4d93-4805-9fc9-43f02af01ecc.md" `
    "c69f6660-4d93-4805-9fc9-43f02af01ecc.7eea7009a8b3a6812f5e9ecd74567f91108cdb53.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1d5adf1f307ce4fa8cd16b204798edccfc83ce99/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/c69f6660-4d93-4805-9fc9-43f02af01ecc.7eea7009a8b3a6812f5e9ecd74567f91108cdb53.zh-cn.xlf" `
    "2016-02-24 09:19:25"

Set-HandbackRow "de-de" `
    "Handed back: in sync with en-US" `
    "c69f6660-4d93-4805-9fc9-43f02af01ecc.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/55e2dd23eb0bd25de4cc8458b6ae18eef6df06dd/e2e/c69f6660-4d93-4805-9fc9-43f02af01ecc.md" `
    "c69f6660-4d93-4805-9fc9-43f02af01ecc.7eea7009a8b3a6812f5e9ecd74567f91108cdb53.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e81e90b747347c8ff53c29da099b8c552d0f8fad/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/c69f6660-4d93-4805-9fc9-43f02af01ecc.7eea7009a8b3a6812f5e9ecd74567f91108cdb53.de-de.xlf" `
    "2016-02-24 09:19:53"

Write-Host "Handback report generated."
